$d = $word.ActiveDocument

# Locate the "Prueba 6" heading paragraph (style "Ttulo2") - this starts the
# block of paragraphs that must be removed (the whole "Prueba 6" test case,
# through the trailing empty paragraph right before the section break).
$startPar = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Prueba 6") {
        $startPar = $p
        break
    }
}

if ($startPar -ne $null) {
    $count = $d.Paragraphs.Count
    $endPar = $d.Paragraphs.Item($count)

    $r = $d.Range($startPar.Range.Start, $endPar.Range.End)
    $r.Delete()
}
